$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) to "UniformF"
$ws.Name = "UniformF"

# Add a new row 16 that duplicates row 15's content/formatting, except
# column A increments to 14 (row 15 has 13), reusing the same text label
# in column B ("HexGrid-60degTilt5degRes").
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
